$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ValidLoginData")
$ws2 = $wb.Worksheets.Item("InvalidLoginData")

$ws2.Range("A5").Value = "shiva"
$ws2.Range("B5").Value = "shiva123"
$ws2.Range("A6").Value = "sri"
$ws2.Range("B6").Value = "sri123"

$ws1.Range("B2").Select()
$ws2.Activate()
